# Update Battery_Data sheet (sheet1) values in column B, rows 2-5
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Battery_Data")

$ws1.Range("B2").Value = 961.4826571789999
$ws1.Range("B3").Value = 711.49716631246
$ws1.Range("B4").Value = 14.2299433262492
$ws1.Range("B5").Value = 118.668407446

# Update Yearly BRC sheet (sheet2): change first five values and delete rows 7-21
$ws2 = $wb.Worksheets.Item("Yearly BRC")

$ws2.Range("B2").Value = 25.31201466157464
$ws2.Range("B3").Value = 25.87534405914255
$ws2.Range("B4").Value = 26.44599501734016
$ws2.Range("B5").Value = 27.03532082871376
$ws2.Range("B6").Value = 27.70666213424036

# Remove rows 7 through 21 (previously y = 6 through y = 20)
$ws2.Range("A7:B21").EntireRow.Delete()
